# Apply scheduled market-data refresh to Golem Profits leve sheets.
# Generated from OOXML diff: updates currentAveragePrice* (H-L) and
# LeveProfit* (M/N) columns; some rows gain/lose an M or N cell
# depending on whether an HQ/NQ price tier is present.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 111
$ws.Range("H111").Value = 1299.6666
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 1299.6666
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 3898.9998
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -10032.9998
# Row 132
$ws.Range("H132").Value = 3979.4
$ws.Range("I132").Value = 1299.6666
$ws.Range("K132").Value = 3898.9998
$ws.Range("M132").Value = -1368.9998
# Row 135
$ws.Range("H135").Value = 536.8889
$ws.Range("I135").Value = 319.16666
$ws.Range("J135").Value = 645.75
$ws.Range("K135").Value = 2872.49994
$ws.Range("L135").Value = 5811.75
$ws.Range("M135").Value = -337.4999399999997
$ws.Range("N135").Value = -10881.75
# Row 137
$ws.Range("H137").Value = 3274.111
$ws.Range("I137").Value = 3479.4
$ws.Range("J137").Value = 3017.5
$ws.Range("K137").Value = 10438.2
$ws.Range("L137").Value = 9052.5
$ws.Range("M137").Value = -7888.200000000001
$ws.Range("N137").Value = -14152.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 781.3333
$ws.Range("I45").Value = 781.3333
$ws.Range("K45").Value = 781.3333
$ws.Range("M45").Value = -404.3333
# Row 61
$ws.Range("H61").Value = 1186
$ws.Range("I61").Value = 1186
$ws.Range("K61").Value = 1186
$ws.Range("M61").Value = -974
# Row 97
$ws.Range("H97").Value = 1722.3
$ws.Range("I97").Value = 1347
$ws.Range("K97").Value = 1347
$ws.Range("M97").Value = -851
# Row 110
$ws.Range("H110").Value = 703.4286
$ws.Range("I110").Value = 703.4286
$ws.Range("K110").Value = 703.4286
$ws.Range("M110").Value = 1341.5714
# Row 136
$ws.Range("H136").Value = 1186
$ws.Range("I136").Value = 1186
$ws.Range("K136").Value = 3558
$ws.Range("M136").Value = -1008

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2400
$ws.Range("J94").Value = 3000
$ws.Range("L94").Value = 3000
$ws.Range("N94").Value = -3902
# Row 99
$ws.Range("H99").Value = 2755.5
$ws.Range("I99").Value = 2606
$ws.Range("K99").Value = 2606
$ws.Range("M99").Value = -1108
# Row 134
$ws.Range("H134").Value = 2908.5
$ws.Range("I134").Value = 2752.5715
$ws.Range("K134").Value = 8257.7145
$ws.Range("M134").Value = -5722.7145

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 971.8570999999999
$ws.Range("I16").Value = 1147.5
$ws.Range("J16").Value = 737.6667
$ws.Range("K16").Value = 1147.5
$ws.Range("L16").Value = 737.6667
$ws.Range("M16").Value = -860.5
$ws.Range("N16").Value = -1311.6667
# Row 31
$ws.Range("H31").Value = 10159.444
$ws.Range("I31").Value = 5846
$ws.Range("J31").Value = 13124.9375
$ws.Range("K31").Value = 5846
$ws.Range("L31").Value = 13124.9375
$ws.Range("M31").Value = -5551
$ws.Range("N31").Value = -13714.9375
# Row 33
$ws.Range("H33").Value = 8224.777
$ws.Range("I33").Value = 1017.8
$ws.Range("J33").Value = 17233.5
$ws.Range("K33").Value = 1017.8
$ws.Range("L33").Value = 17233.5
$ws.Range("M33").Value = -638.8
$ws.Range("N33").Value = -17991.5
# Row 34
$ws.Range("H34").Value = 10159.444
$ws.Range("I34").Value = 5846
$ws.Range("J34").Value = 13124.9375
$ws.Range("K34").Value = 5846
$ws.Range("L34").Value = 13124.9375
$ws.Range("M34").Value = -5644
$ws.Range("N34").Value = -13528.9375
# Row 42
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
# Row 44
$ws.Range("H44").Value = 30000
$ws.Range("I44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("M44").Value = -29558
$ws.Range("N44").Value = -30884
# Row 86
$ws.Range("H86").Value = 5849.4
$ws.Range("J86").Value = 5474.5
$ws.Range("L86").Value = 5474.5
$ws.Range("N86").Value = -7720.5
# Row 89
$ws.Range("H89").Value = 5849.4
$ws.Range("J89").Value = 5474.5
$ws.Range("L89").Value = 27372.5
$ws.Range("N89").Value = -38604.5
# Row 113
$ws.Range("H113").Value = 971.8570999999999
$ws.Range("I113").Value = 1147.5
$ws.Range("J113").Value = 737.6667
$ws.Range("K113").Value = 1147.5
$ws.Range("L113").Value = 737.6667
$ws.Range("M113").Value = 1022.5
$ws.Range("N113").Value = -5077.6667
# Row 132
$ws.Range("H132").Value = 1684
$ws.Range("I132").Value = 1605
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4815
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2285
$ws.Range("N132").Value = -11060
# Row 134
$ws.Range("H134").Value = 1020.85
$ws.Range("I134").Value = 959.8823
$ws.Range("K134").Value = 2879.6469
$ws.Range("M134").Value = -344.6468999999997

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 41109.7
$ws.Range("I102").Value = 45121.89
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 45121.89
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -43499.89
$ws.Range("N102").Value = -8244
# Row 132
$ws.Range("H132").Value = 1804
$ws.Range("I132").Value = 1804
$ws.Range("K132").Value = 5412
$ws.Range("M132").Value = -2882

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2116.25
$ws.Range("I22").Value = 1416.5
$ws.Range("J22").Value = 2816
$ws.Range("K22").Value = 1416.5
$ws.Range("L22").Value = 2816
$ws.Range("M22").Value = -1121.5
$ws.Range("N22").Value = -3406
# Row 27
$ws.Range("H27").Value = 2116.25
$ws.Range("I27").Value = 1416.5
$ws.Range("J27").Value = 2816
$ws.Range("K27").Value = 1416.5
$ws.Range("L27").Value = 2816
$ws.Range("M27").Value = -1309.5
$ws.Range("N27").Value = -3030
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 32
$ws.Range("H32").Value = 2026
$ws.Range("I32").Value = 2026
$ws.Range("K32").Value = 2026
$ws.Range("M32").Value = -1709
# Row 64
$ws.Range("H64").Value = 19341.75
$ws.Range("I64").Value = 30000
$ws.Range("K64").Value = 30000
$ws.Range("M64").Value = -29752
# Row 67
$ws.Range("H67").Value = 19341.75
$ws.Range("I67").Value = 30000
$ws.Range("K67").Value = 30000
$ws.Range("M67").Value = -29142
# Row 132
$ws.Range("H132").Value = 1604.3
$ws.Range("I132").Value = 1630.5
$ws.Range("K132").Value = 4891.5
$ws.Range("M132").Value = -2361.5
